# Update crypto price/volume table (rows 2-51) to reflect the latest
# scrape from the GitHub Actions job. Column D ("Price") values that
# look like plain decimal numbers are prefixed with a literal leading
# apostrophe so Excel stores them as quote-prefixed TEXT (matching the
# source workbook, which keeps every price as a literal string such as
# "1.000" or "0.00000000120") instead of silently coercing them to
# numbers and losing the exact formatting/trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.246.76'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '1.869.86'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("D4").Value = '''0.9999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''0.7104'
$ws.Range("E5").Value = '  -1.80%  '
$ws.Range("D6").Value = '''241.50'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '''0.3103'
$ws.Range("E8").Value = '  +0.50%  '
$ws.Range("D9").Value = '''0.07689'
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").Value = '''25.00'
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("D11").Value = '''0.08383'
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").Value = '1.883.72'
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").Value = '''5.203'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '''0.7096'
$ws.Range("E14").Value = '  -1.66%  '
$ws.Range("D15").Value = '''91.12'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '29.256.02'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '''0.000008284'
$ws.Range("E17").Value = '  +6.08%  '
$ws.Range("D18").Value = '''5.931'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("D19").Value = '''242.14'
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").Value = '2.129.64'
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("D21").Value = '''13.16'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").Value = '''0.9998'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '''7.821'
$ws.Range("E23").Value = '  -1.82%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").Value = '''0.1635'
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("D26").Value = '''163.08'
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("D27").Value = '''8.998'
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("D28").Value = '''18.47'
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("D30").Value = '''4.404'
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = '''4.322'
$ws.Range("E31").Value = '  +5.20%  '
$ws.Range("D32").Value = '''1.282'
$ws.Range("E32").Value = '  -4.54%  '
$ws.Range("D33").Value = '''0.05234'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").Value = '''1.922'
$ws.Range("D35").Value = '''0.7514'
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("D36").Value = '''1.170'
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").Value = '''2.683'
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").Value = '1.154.70'
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("D41").Value = '''6.356'
$ws.Range("E41").Value = '  +3.81%  '
$ws.Range("D42").Value = '''73.04'
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("D43").Value = '''0.8861'
$ws.Range("E43").Value = '  -1.93%  '
$ws.Range("D44").Value = '''104.04'
$ws.Range("E44").Value = '  +1.97%  '
$ws.Range("D45").Value = '''0.9997'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").Value = '2.026.50'
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.793'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '''0.5185'
$ws.Range("E48").Value = '  -1.91%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000120'
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.377'
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '''0.4290'
$ws.Range("E51").Value = '  +0.58%  '
